# stock_data.xlsx — add "revenue" and "people" columns, rewrite the pe->revenue
# column, and append two more rows (RIL / TATA) per the "added read write excel
# file folder" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row --------------------------------------------------------
$ws.Range("A1").Value = "tickers"
$ws.Range("B1").Value = "eps"
$ws.Range("C1").Value = "revenue"
$ws.Range("D1").Value = "price"
$ws.Range("E1").Value = "people"

# ---- Data rows ----------------------------------------------------------
$ws.Range("A2").Value = "GOOGL"
$ws.Range("B2").Value = 27.82
$ws.Range("C2").Value = 87
$ws.Range("D2").Value = 845
$ws.Range("E2").Value = "larry page"

$ws.Range("A3").Value = "WMT"
$ws.Range("B3").Value = 4.61
$ws.Range("C3").Value = 484
$ws.Range("D3").Value = 65
$ws.Range("E3").Value = "n.a."

$ws.Range("A4").Value = "MSFT"
$ws.Range("B4").Value = -1
$ws.Range("C4").Value = 85
$ws.Range("D4").Value = 64
$ws.Range("E4").Value = "bill gates"

$ws.Range("A5").Value = "RIL "
$ws.Range("B5").Value = "not available"
$ws.Range("C5").Value = 50
$ws.Range("D5").Value = 1023
$ws.Range("E5").Value = "mukesh ambani"

$ws.Range("A6").Value = "TATA"
$ws.Range("B6").Value = 5.6
$ws.Range("C6").Value = -1
$ws.Range("D6").Value = "n.a."
$ws.Range("E6").Value = "ratan tata"

# ---- Header formatting (bold, matches the existing A1:D1 style) ---------
$ws.Range("A1:E1").Font.Bold = $true

# ---- Column widths (best-fit on the two text-heavy columns) -------------
$ws.Columns.Item(2).ColumnWidth = 11.584
$ws.Columns.Item(5).ColumnWidth = 14.25

# ---- View state: zoom + active selection --------------------------------
$win = $ws.Application.ActiveWindow
$win.Zoom = 160
$ws.Range("D12").Select()
